# =============================================================================
# voltage commutation possible with ecu v1.0
#  - rename Sheet1 -> ADC, add a new PWM sheet
#  - ADC: add per-channel timing helper formulas (C2, C4, D4), tweak the
#    sample-time / channel-count inputs, restyle a few numeric outputs
#  - PWM: brand-new timer/prescaler worksheet
# =============================================================================

$wb  = $excel.ActiveWorkbook

$adc = $wb.Worksheets.Item(1)
$adc.Name = "ADC"

$pwm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $adc)
$pwm.Name = "PWM"

# =============================================================================
# ADC sheet (formerly "Sheet1")
# =============================================================================

# New helper formulas in columns C/D
$adc.Range("C2").Formula = "=1/B2"
$adc.Range("C2").NumberFormat = "0.000000000"

$adc.Range("C4").Formula = "=1/B4"
$adc.Range("C4").NumberFormat = "0.000000000000"

$adc.Range("D4").Formula = "=C2*4"
$adc.Range("D4").NumberFormat = "0.000000000"

$adc.Range("B4").NumberFormat = "0"

# Updated inputs: sample time per channel (cycles) 3 -> 84, nr channels 7 -> 4
$adc.Range("B7").Value = 84
$adc.Range("B9").Value = 4

# Restyle the total conversion time and drop the stray D11 formula
$adc.Range("B11").NumberFormat = "0.00000000000"
$adc.Range("D11").ClearContents()

# Column widths (C/D newly used)
$adc.Columns.Item(3).ColumnWidth = 17.33
$adc.Columns.Item(4).ColumnWidth = 13

# View: zoom in on the ADC sheet, active cell C2
$adc.Activate()
$excel.ActiveWindow.Zoom = 220
$adc.Range("C2").Select()

# =============================================================================
# PWM sheet (new)
# =============================================================================

$pwm.Range("A1:C1").Merge()
$pwm.Range("A1").Value = "PWM period"
$pwm.Range("A1:C1").HorizontalAlignment = -4108   # xlCenter

$pwm.Range("A2").Value = "APB2 Freq"
$pwm.Range("B2").Value = 84000000

$pwm.Range("A3").Value = "Prescaler"
$pwm.Range("B3").Value = 0

$pwm.Range("A4").Value = "TIM Freq"
$pwm.Range("B4").Formula = "=B2/(B3+1)"

$pwm.Range("G4").NumberFormat = "#,##0.00 ""lei"";[Red]-#,##0.00 ""lei"""

$pwm.Range("A5").Value = "Tcnt"
$pwm.Range("B5").Formula = "=(1/B4)*1000000"
$pwm.Range("B5").NumberFormat = "0.000000000000"

$pwm.Range("D5").Value = "ADC Sampling duration"
$pwm.Range("E5").Formula = "=1290*B5"
$pwm.Range("E5").ClearFormats()

$pwm.Range("B6").NumberFormat = "0.000000000000"
$pwm.Range("B7").NumberFormat = "0.000000000000"
$pwm.Range("B8").NumberFormat = "0.00000"

$pwm.Range("A9").Value = "ARR"
$pwm.Range("B9").Value = 4199

$pwm.Range("A10").Value = "PWM Freq"
$pwm.Range("B10").Formula = "=B4/(B9+1)"

# Column widths
$pwm.Columns.Item(1).ColumnWidth = 9.33
$pwm.Columns.Item(2).ColumnWidth = 18.5
$pwm.Columns.Item(4).ColumnWidth = 21.83
$pwm.Columns.Item(5).ColumnWidth = 11.33

# View: zoom in on the PWM sheet, active cell E6, make it the selected/visible tab
$pwm.Activate()
$excel.ActiveWindow.Zoom = 145
$pwm.Range("E6").Select()
